$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.794.64"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "2.538.89"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.34"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.87"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "2.535.59"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.12"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.54"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "2.918.86"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "67.643.47"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "2.516.16"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.86"
$ws.Range("E19").Value = "  +5.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.89"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.82"
$ws.Range("E21").Value = "  +4.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.60"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.81"
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.98"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("D29").Value = "2.651.99"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "0.0₃0967"
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.44"
$ws.Range("E31").Value = "  +4.53%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "538.68"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.129"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.64"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.10"
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.63"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.352"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.16"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("E46").Value = "  +4.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "147.05"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.554"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.72"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.73"
$ws.Range("E50").Value = "  +2.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0748"
$ws.Range("E51").Value = "  -1.04%  "